$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old header row (row 3: "stt", "Name", "Age", "City").
# This shifts all the data rows below it up by one.
$ws.Rows.Item(3).Delete()

# Row 2 (previously empty) becomes the new header row.
$ws.Range("A2").Value = "ID"
$ws.Range("B2").Value = "Name"
$ws.Range("C2").Value = "Age"
$ws.Range("D2").Value = "City"
